# Update "想去人数" (number of interested attendees) figures on the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, plus a single
# matching correction on 演出 (Shows), per upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 294
$ws1.Range("F7").Value = 1448
$ws1.Range("F8").Value = 589
$ws1.Range("F10").Value = 754
$ws1.Range("F12").Value = 175
$ws1.Range("F15").Value = 1391
$ws1.Range("F18").Value = 281
$ws1.Range("F20").Value = 74
$ws1.Range("F22").Value = 1013
$ws1.Range("F23").Value = 39
$ws1.Range("F24").Value = 243
$ws1.Range("F26").Value = 5988
$ws1.Range("F28").Value = 126
$ws1.Range("F31").Value = 14690
$ws1.Range("F32").Value = 1457
$ws1.Range("F36").Value = 9488
$ws1.Range("F37").Value = 647

# --- Sheet: 演出 (Shows) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 341

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 294
$ws4.Range("F7").Value = 1448
$ws4.Range("F8").Value = 589
$ws4.Range("F10").Value = 754
$ws4.Range("F12").Value = 175
$ws4.Range("F15").Value = 1391
$ws4.Range("F18").Value = 281
$ws4.Range("F19").Value = 341
$ws4.Range("F21").Value = 74
$ws4.Range("F24").Value = 1013
$ws4.Range("F25").Value = 39
$ws4.Range("F26").Value = 243
$ws4.Range("F29").Value = 5988
$ws4.Range("F31").Value = 126
$ws4.Range("F34").Value = 14690
$ws4.Range("F35").Value = 1457
$ws4.Range("F39").Value = 9488
$ws4.Range("F40").Value = 647
